$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 207 (existing rows 207-222 shift down to 208-223).
$ws.Rows.Item(207).Insert()

# Populate the newly-inserted row 207 with the new weekly data point.
$ws.Cells.Item(207, 1).Value = 4
$ws.Cells.Item(207, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(207, 3).Value = "Los Lagos"
$ws.Cells.Item(207, 4).Value = 44585
$ws.Cells.Item(207, 5).Value = 10
$ws.Cells.Item(207, 6).Value = 100112040
$ws.Cells.Item(207, 7).Value = "Cilantro"
$ws.Cells.Item(207, 8).Value = "Sin especificar"
$ws.Cells.Item(207, 9).Value = "Primera"
$ws.Cells.Item(207, 10).Value = 60
$ws.Cells.Item(207, 11).Value = 10000
$ws.Cells.Item(207, 12).Value = 10000
$ws.Cells.Item(207, 13).Value = 10000
$ws.Cells.Item(207, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(207, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(207, 16).Value = 5000
$ws.Cells.Item(207, 17).Value = 2
$ws.Cells.Item(207, 18).Value = "Hortaliza"
